$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(28, 8).Value = 713.26666
$ws.Cells.Item(28, 9).Value = 671.44446
$ws.Cells.Item(28, 10).Value = 776
$ws.Cells.Item(28, 11).Value = 671.44446
$ws.Cells.Item(28, 12).Value = 776
$ws.Cells.Item(28, 13).Value = -186.44446
$ws.Cells.Item(28, 14).Value = -1746

$ws.Cells.Item(42, 8).Value = 216.42857
$ws.Cells.Item(42, 9).Value = 201
$ws.Cells.Item(42, 10).Value = 309
$ws.Cells.Item(42, 11).Value = 603
$ws.Cells.Item(42, 12).Value = 927
$ws.Cells.Item(42, 13).Value = -373
$ws.Cells.Item(42, 14).Value = -1387

$ws.Cells.Item(51, 8).Value = 4307
$ws.Cells.Item(51, 9).Value = 1933.5
$ws.Cells.Item(51, 11).Value = 1933.5
$ws.Cells.Item(51, 13).Value = -1449.5

$ws.Cells.Item(64, 8).Value = 4106.564
$ws.Cells.Item(64, 9).Value = 4188.5713
$ws.Cells.Item(64, 10).Value = 3897.818
$ws.Cells.Item(64, 11).Value = 4188.5713
$ws.Cells.Item(64, 12).Value = 3897.818
$ws.Cells.Item(64, 13).Value = -3940.5713
$ws.Cells.Item(64, 14).Value = -4393.818

$ws.Cells.Item(67, 8).Value = 4106.564
$ws.Cells.Item(67, 9).Value = 4188.5713
$ws.Cells.Item(67, 10).Value = 3897.818
$ws.Cells.Item(67, 11).Value = 4188.5713
$ws.Cells.Item(67, 12).Value = 3897.818
$ws.Cells.Item(67, 13).Value = -3330.5713
$ws.Cells.Item(67, 14).Value = -5613.818

$ws.Cells.Item(86, 8).Value = 4399.1816
$ws.Cells.Item(86, 9).Value = 2751.375
$ws.Cells.Item(86, 10).Value = 8793.333000000001
$ws.Cells.Item(86, 11).Value = 2751.375
$ws.Cells.Item(86, 12).Value = 8793.333000000001
$ws.Cells.Item(86, 13).Value = -1628.375
$ws.Cells.Item(86, 14).Value = -11039.333

$ws.Cells.Item(89, 8).Value = 4399.1816
$ws.Cells.Item(89, 9).Value = 2751.375
$ws.Cells.Item(89, 10).Value = 8793.333000000001
$ws.Cells.Item(89, 11).Value = 13756.875
$ws.Cells.Item(89, 12).Value = 43966.665
$ws.Cells.Item(89, 13).Value = -8140.875
$ws.Cells.Item(89, 14).Value = -55198.665

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(47, 8).Value = 24900
$ws.Cells.Item(47, 10).Value = 24900
$ws.Cells.Item(47, 12).Value = 24900
$ws.Cells.Item(47, 14).Value = -26350

$ws.Cells.Item(63, 8).Value = 166668100
$ws.Cells.Item(63, 9).Value = 166668100
$ws.Cells.Item(63, 10).Value = 0
$ws.Cells.Item(63, 11).Value = 166668100
$ws.Cells.Item(63, 12).Value = 0
$ws.Cells.Item(63, 13).Value = -166667414
$ws.Cells.Item(63, 14).Value = ""

$ws.Cells.Item(66, 8).Value = 166668100
$ws.Cells.Item(66, 9).Value = 166668100
$ws.Cells.Item(66, 10).Value = 0
$ws.Cells.Item(66, 11).Value = 833340500
$ws.Cells.Item(66, 12).Value = 0
$ws.Cells.Item(66, 13).Value = -833337068
$ws.Cells.Item(66, 14).Value = ""

$ws.Cells.Item(74, 8).Value = 1573.8055
$ws.Cells.Item(74, 9).Value = 1717.3077
$ws.Cells.Item(74, 10).Value = 1492.6957
$ws.Cells.Item(74, 11).Value = 1717.3077
$ws.Cells.Item(74, 12).Value = 1492.6957
$ws.Cells.Item(74, 13).Value = -843.3077000000001
$ws.Cells.Item(74, 14).Value = -3240.6957

$ws.Cells.Item(77, 8).Value = 1573.8055
$ws.Cells.Item(77, 9).Value = 1717.3077
$ws.Cells.Item(77, 10).Value = 1492.6957
$ws.Cells.Item(77, 11).Value = 8586.538500000001
$ws.Cells.Item(77, 12).Value = 7463.4785
$ws.Cells.Item(77, 13).Value = -4218.538500000001
$ws.Cells.Item(77, 14).Value = -16199.4785

$ws.Cells.Item(88, 8).Value = 100002456
$ws.Cells.Item(88, 9).Value = 2747.2
$ws.Cells.Item(88, 10).Value = 200002160
$ws.Cells.Item(88, 11).Value = 2747.2
$ws.Cells.Item(88, 12).Value = 200002160
$ws.Cells.Item(88, 13).Value = -2341.2
$ws.Cells.Item(88, 14).Value = -200002972

$ws.Cells.Item(91, 8).Value = 100002456
$ws.Cells.Item(91, 9).Value = 2747.2
$ws.Cells.Item(91, 10).Value = 200002160
$ws.Cells.Item(91, 11).Value = 2747.2
$ws.Cells.Item(91, 12).Value = 200002160
$ws.Cells.Item(91, 13).Value = -1343.2
$ws.Cells.Item(91, 14).Value = -200004968

$ws.Cells.Item(132, 8).Value = 2801.5686
$ws.Cells.Item(132, 9).Value = 2420.8235
$ws.Cells.Item(132, 10).Value = 3563.0588
$ws.Cells.Item(132, 11).Value = 7262.470499999999
$ws.Cells.Item(132, 12).Value = 10689.1764
$ws.Cells.Item(132, 13).Value = -4732.470499999999
$ws.Cells.Item(132, 14).Value = -15749.1764

$ws.Cells.Item(138, 8).Value = 44426.668
$ws.Cells.Item(138, 10).Value = 44426.668
$ws.Cells.Item(138, 12).Value = 44426.668
$ws.Cells.Item(138, 14).Value = -54706.668

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(86, 8).Value = 1970.8
$ws.Cells.Item(86, 9).Value = 1795.75
$ws.Cells.Item(86, 10).Value = 2671
$ws.Cells.Item(86, 11).Value = 1795.75
$ws.Cells.Item(86, 12).Value = 2671
$ws.Cells.Item(86, 13).Value = -672.75
$ws.Cells.Item(86, 14).Value = -4917

$ws.Cells.Item(89, 8).Value = 1970.8
$ws.Cells.Item(89, 9).Value = 1795.75
$ws.Cells.Item(89, 10).Value = 2671
$ws.Cells.Item(89, 11).Value = 8978.75
$ws.Cells.Item(89, 12).Value = 13355
$ws.Cells.Item(89, 13).Value = -3362.75
$ws.Cells.Item(89, 14).Value = -24587

$ws.Cells.Item(134, 8).Value = 3626
$ws.Cells.Item(134, 9).Value = 5504
$ws.Cells.Item(134, 10).Value = 3000
$ws.Cells.Item(134, 11).Value = 16512
$ws.Cells.Item(134, 12).Value = 9000
$ws.Cells.Item(134, 13).Value = -13977
$ws.Cells.Item(134, 14).Value = -14070

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 3751.25
$ws.Cells.Item(31, 10).Value = 4536.364
$ws.Cells.Item(31, 12).Value = 4536.364
$ws.Cells.Item(31, 14).Value = -5126.364

$ws.Cells.Item(34, 8).Value = 3751.25
$ws.Cells.Item(34, 10).Value = 4536.364
$ws.Cells.Item(34, 12).Value = 4536.364
$ws.Cells.Item(34, 14).Value = -4940.364

$ws.Cells.Item(62, 8).Value = 3775.6667
$ws.Cells.Item(62, 9).Value = 3274.75
$ws.Cells.Item(62, 10).Value = 4777.5
$ws.Cells.Item(62, 11).Value = 3274.75
$ws.Cells.Item(62, 12).Value = 4777.5
$ws.Cells.Item(62, 13).Value = -2650.75
$ws.Cells.Item(62, 14).Value = -6025.5

$ws.Cells.Item(63, 8).Value = 40100
$ws.Cells.Item(63, 10).Value = 40100
$ws.Cells.Item(63, 12).Value = 40100
$ws.Cells.Item(63, 14).Value = -41472

$ws.Cells.Item(65, 8).Value = 3775.6667
$ws.Cells.Item(65, 9).Value = 3274.75
$ws.Cells.Item(65, 10).Value = 4777.5
$ws.Cells.Item(65, 11).Value = 16373.75
$ws.Cells.Item(65, 12).Value = 23887.5
$ws.Cells.Item(65, 13).Value = -13253.75
$ws.Cells.Item(65, 14).Value = -30127.5

$ws.Cells.Item(66, 8).Value = 40100
$ws.Cells.Item(66, 10).Value = 40100
$ws.Cells.Item(66, 12).Value = 120300
$ws.Cells.Item(66, 14).Value = -127164

$ws.Cells.Item(99, 8).Value = 24445.363
$ws.Cells.Item(99, 9).Value = 22000
$ws.Cells.Item(99, 10).Value = 26483.166
$ws.Cells.Item(99, 11).Value = 22000
$ws.Cells.Item(99, 12).Value = 26483.166
$ws.Cells.Item(99, 13).Value = -20502
$ws.Cells.Item(99, 14).Value = -29479.166

$ws.Cells.Item(126, 8).Value = 24445.363
$ws.Cells.Item(126, 9).Value = 22000
$ws.Cells.Item(126, 10).Value = 26483.166
$ws.Cells.Item(126, 11).Value = 66000
$ws.Cells.Item(126, 12).Value = 79449.49800000001
$ws.Cells.Item(126, 13).Value = -63530
$ws.Cells.Item(126, 14).Value = -84389.49800000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(97, 8).Value = 20000340
$ws.Cells.Item(97, 9).Value = 20000340
$ws.Cells.Item(97, 11).Value = 60001020
$ws.Cells.Item(97, 13).Value = -60000524

$ws.Cells.Item(107, 8).Value = 680.25
$ws.Cells.Item(107, 9).Value = 216.33333
$ws.Cells.Item(107, 10).Value = 1144.1666
$ws.Cells.Item(107, 11).Value = 648.99999
$ws.Cells.Item(107, 12).Value = 3432.4998
$ws.Cells.Item(107, 13).Value = 1271.00001
$ws.Cells.Item(107, 14).Value = -7272.4998

$ws.Cells.Item(131, 8).Value = 23914128
$ws.Cells.Item(131, 10).Value = 27779040
$ws.Cells.Item(131, 12).Value = 83337120
$ws.Cells.Item(131, 14).Value = -83347200

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(80, 8).Value = 2581.081
$ws.Cells.Item(80, 9).Value = 2582.6086
$ws.Cells.Item(80, 11).Value = 2582.6086
$ws.Cells.Item(80, 13).Value = -1584.6086

$ws.Cells.Item(83, 8).Value = 2581.081
$ws.Cells.Item(83, 9).Value = 2582.6086
$ws.Cells.Item(83, 11).Value = 12913.043
$ws.Cells.Item(83, 13).Value = -7921.043

$ws.Cells.Item(122, 8).Value = 62638044
$ws.Cells.Item(122, 10).Value = 2401.6
$ws.Cells.Item(122, 12).Value = 7204.799999999999
$ws.Cells.Item(122, 14).Value = -12104.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(40, 8).Value = 142860290
$ws.Cells.Item(40, 9).Value = 250002750
$ws.Cells.Item(40, 10).Value = 3668.3333
$ws.Cells.Item(40, 11).Value = 250002750
$ws.Cells.Item(40, 12).Value = 3668.3333
$ws.Cells.Item(40, 13).Value = -250002614
$ws.Cells.Item(40, 14).Value = -3940.3333

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(122, 8).Value = 1523.7333
$ws.Cells.Item(122, 9).Value = 1315.1
$ws.Cells.Item(122, 10).Value = 1941
$ws.Cells.Item(122, 11).Value = 3945.3
$ws.Cells.Item(122, 12).Value = 5823
$ws.Cells.Item(122, 13).Value = -1495.3
$ws.Cells.Item(122, 14).Value = -10723

$ws.Cells.Item(126, 8).Value = 1188.6666
$ws.Cells.Item(126, 9).Value = 782.3333
$ws.Cells.Item(126, 11).Value = 2346.9999
$ws.Cells.Item(126, 13).Value = 123.0001000000002
